$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95 <= (old) Row 96 data  [id 6782567]
$ws.Range('B95').Value = 6782567
$ws.Range('E95').Value = 'AD Grecia'
$ws.Range('F95').Value = 'Municipal Liberia'
$ws.Range('G95').Value = 2
$ws.Range('H95').Value = 3
$ws.Range('I95').Value = 'A'
$ws.Range('J95').Value = 2.875
$ws.Range('K95').Value = 3.5
$ws.Range('L95').Value = 2.15
$ws.Range('M95').Value = 2.3
$ws.Range('N95').Value = 3.5
$ws.Range('O95').Value = 2.6
$ws.Range('P95').Value = 0
$ws.Range('Q95').Value = 1.8
$ws.Range('R95').Value = 2
$ws.Range('S95').Value = 2.75
$ws.Range('T95').Value = 1.8
$ws.Range('U95').Value = 2
$ws.Range('V95').Value = -1
$ws.Range('W95').Value = -1
$ws.Range('X95').Value = 1.6
$ws.Range('Y95').Value = -1
$ws.Range('Z95').Value = 1
$ws.Range('AA95').Value = 0.8
$ws.Range('AB95').Value = -1

# Row 96 <= (old) Row 95 data  [id 6782565]
$ws.Range('B96').Value = 6782565
$ws.Range('E96').Value = 'Santos de Gupiles'
$ws.Range('F96').Value = 'Municipal Perez Zeledon'
$ws.Range('G96').Value = 2
$ws.Range('H96').Value = 0
$ws.Range('I96').Value = 'H'
$ws.Range('J96').Value = 1.833
$ws.Range('K96').Value = 3.4
$ws.Range('L96').Value = 3.75
$ws.Range('M96').Value = 1.833
$ws.Range('N96').Value = 3.5
$ws.Range('O96').Value = 3.75
$ws.Range('P96').Value = -0.5
$ws.Range('Q96').Value = 1.875
$ws.Range('R96').Value = 1.925
$ws.Range('S96').Value = 2.75
$ws.Range('T96').Value = 2
$ws.Range('U96').Value = 1.8
$ws.Range('V96').Value = 0.833
$ws.Range('W96').Value = -1
$ws.Range('X96').Value = -1
$ws.Range('Y96').Value = 0.875
$ws.Range('Z96').Value = -1
$ws.Range('AA96').Value = -1
$ws.Range('AB96').Value = 0.8

# Row 110 <= (old) Row 111 data  [id 6782579]
$ws.Range('B110').Value = 6782579
$ws.Range('E110').Value = 'Santos de Gupiles'
$ws.Range('F110').Value = 'AD San Carlos'
$ws.Range('G110').Value = 0
$ws.Range('H110').Value = 2
$ws.Range('I110').Value = 'A'
$ws.Range('J110').Value = 2.4
$ws.Range('K110').Value = 3.3
$ws.Range('L110').Value = 2.7
$ws.Range('M110').Value = 2.375
$ws.Range('N110').Value = 3.4
$ws.Range('O110').Value = 2.8
$ws.Range('P110').Value = -0.25
$ws.Range('Q110').Value = 2
$ws.Range('R110').Value = 1.8
$ws.Range('S110').Value = 2.5
$ws.Range('T110').Value = 1.875
$ws.Range('U110').Value = 1.925
$ws.Range('V110').Value = -1
$ws.Range('W110').Value = -1
$ws.Range('X110').Value = 1.8
$ws.Range('Y110').Value = -1
$ws.Range('Z110').Value = 0.8
$ws.Range('AA110').Value = -1
$ws.Range('AB110').Value = 0.925

# Row 111 <= (old) Row 110 data  [id 6782581]
$ws.Range('B111').Value = 6782581
$ws.Range('E111').Value = 'Alajuelense'
$ws.Range('F111').Value = 'AD Grecia'
$ws.Range('G111').Value = 2
$ws.Range('H111').Value = 0
$ws.Range('I111').Value = 'H'
$ws.Range('J111').Value = 1.181
$ws.Range('K111').Value = 6.5
$ws.Range('L111').Value = 11
$ws.Range('M111').Value = 1.25
$ws.Range('N111').Value = 5
$ws.Range('O111').Value = 9
$ws.Range('P111').Value = -1.75
$ws.Range('Q111').Value = 1.975
$ws.Range('R111').Value = 1.825
$ws.Range('S111').Value = 3.25
$ws.Range('T111').Value = 2
$ws.Range('U111').Value = 1.8
$ws.Range('V111').Value = 0.25
$ws.Range('W111').Value = -1
$ws.Range('X111').Value = -1
$ws.Range('Y111').Value = 0.4875
$ws.Range('Z111').Value = -0.5
$ws.Range('AA111').Value = -1
$ws.Range('AB111').Value = 0.8

# Row 129 <= (old) Row 131 data  [id 6782596]
$ws.Range('B129').Value = 6782596
$ws.Range('E129').Value = 'Alajuelense'
$ws.Range('F129').Value = 'AD Guanacasteca'
$ws.Range('G129').Value = 3
$ws.Range('H129').Value = 4
$ws.Range('I129').Value = 'A'
$ws.Range('J129').Value = 1.363
$ws.Range('K129').Value = 4.75
$ws.Range('L129').Value = 8
$ws.Range('M129').Value = 1.444
$ws.Range('N129').Value = 4.333
$ws.Range('O129').Value = 7
$ws.Range('P129').Value = -1.25
$ws.Range('Q129').Value = 1.975
$ws.Range('R129').Value = 1.825
$ws.Range('S129').Value = 2.75
$ws.Range('T129').Value = 1.775
$ws.Range('U129').Value = 2.025
$ws.Range('V129').Value = -1
$ws.Range('W129').Value = -1
$ws.Range('X129').Value = 6
$ws.Range('Y129').Value = -1
$ws.Range('Z129').Value = 0.825
$ws.Range('AA129').Value = 0.7749999999999999
$ws.Range('AB129').Value = -1

# Row 131 <= (old) Row 129 data  [id 6782595]
$ws.Range('B131').Value = 6782595
$ws.Range('E131').Value = 'Herediano'
$ws.Range('F131').Value = 'Sporting San Jose'
$ws.Range('G131').Value = 3
$ws.Range('H131').Value = 0
$ws.Range('I131').Value = 'H'
$ws.Range('J131').Value = 1.4
$ws.Range('K131').Value = 4.75
$ws.Range('L131').Value = 7
$ws.Range('M131').Value = 1.363
$ws.Range('N131').Value = 4.75
$ws.Range('O131').Value = 8.5
$ws.Range('P131').Value = -1.25
$ws.Range('Q131').Value = 1.8
$ws.Range('R131').Value = 2
$ws.Range('S131').Value = 3
$ws.Range('T131').Value = 1.95
$ws.Range('U131').Value = 1.85
$ws.Range('V131').Value = 0.363
$ws.Range('W131').Value = -1
$ws.Range('X131').Value = -1
$ws.Range('Y131').Value = 0.8
$ws.Range('Z131').Value = -1
$ws.Range('AA131').Value = 0
$ws.Range('AB131').Value = 0

# Row 192 <= (old) Row 193 data  [id 7623919]
$ws.Range('B192').Value = 7623919
$ws.Range('E192').Value = 'Municipal Liberia'
$ws.Range('F192').Value = 'Sporting San Jose'
$ws.Range('G192').Value = 2
$ws.Range('H192').Value = 0
$ws.Range('I192').Value = 'H'
$ws.Range('J192').Value = 1.75
$ws.Range('K192').Value = 3.6
$ws.Range('L192').Value = 3.8
$ws.Range('M192').Value = 1.8
$ws.Range('N192').Value = 3.6
$ws.Range('O192').Value = 3.6
$ws.Range('P192').Value = -0.5
$ws.Range('Q192').Value = 1.9
$ws.Range('R192').Value = 1.9
$ws.Range('S192').Value = 2.75
$ws.Range('T192').Value = 2
$ws.Range('U192').Value = 1.8
$ws.Range('V192').Value = 0.8
$ws.Range('W192').Value = -1
$ws.Range('X192').Value = -1
$ws.Range('Y192').Value = 0.8999999999999999
$ws.Range('Z192').Value = -1
$ws.Range('AA192').Value = -1
$ws.Range('AB192').Value = 0.8

# Row 193 <= (old) Row 192 data  [id 7623916]
$ws.Range('B193').Value = 7623916
$ws.Range('E193').Value = 'Santos de Gupiles'
$ws.Range('F193').Value = 'AD Grecia'
$ws.Range('G193').Value = 0
$ws.Range('H193').Value = 2
$ws.Range('I193').Value = 'A'
$ws.Range('J193').Value = 2.05
$ws.Range('K193').Value = 3.3
$ws.Range('L193').Value = 3.2
$ws.Range('M193').Value = 1.909
$ws.Range('N193').Value = 3.4
$ws.Range('O193').Value = 3.6
$ws.Range('P193').Value = -0.5
$ws.Range('Q193').Value = 1.95
$ws.Range('R193').Value = 1.85
$ws.Range('S193').Value = 2.5
$ws.Range('T193').Value = 1.85
$ws.Range('U193').Value = 1.95
$ws.Range('V193').Value = -1
$ws.Range('W193').Value = -1
$ws.Range('X193').Value = 2.6
$ws.Range('Y193').Value = -1
$ws.Range('Z193').Value = 0.8500000000000001
$ws.Range('AA193').Value = -1
$ws.Range('AB193').Value = 0.95

# Row 200 <= (old) Row 201 data  [id 7623921]
$ws.Range('B200').Value = 7623921
$ws.Range('E200').Value = 'AD Grecia'
$ws.Range('F200').Value = 'Municipal Liberia'
$ws.Range('G200').Value = 1
$ws.Range('H200').Value = 2
$ws.Range('I200').Value = 'A'
$ws.Range('J200').Value = 2.75
$ws.Range('K200').Value = 3.25
$ws.Range('L200').Value = 2.3
$ws.Range('M200').Value = 3.1
$ws.Range('N200').Value = 3.25
$ws.Range('O200').Value = 2.1
$ws.Range('P200').Value = 0.25
$ws.Range('Q200').Value = 1.9
$ws.Range('R200').Value = 1.9
$ws.Range('S200').Value = 2.5
$ws.Range('T200').Value = 1.9
$ws.Range('U200').Value = 1.9
$ws.Range('V200').Value = -1
$ws.Range('W200').Value = -1
$ws.Range('X200').Value = 1.1
$ws.Range('Y200').Value = -1
$ws.Range('Z200').Value = 0.8999999999999999
$ws.Range('AA200').Value = 0.8999999999999999
$ws.Range('AB200').Value = -1

# Row 201 <= (old) Row 200 data  [id 7624967]
$ws.Range('B201').Value = 7624967
$ws.Range('E201').Value = 'Puntarenas'
$ws.Range('F201').Value = 'Herediano'
$ws.Range('G201').Value = 0
$ws.Range('H201').Value = 0
$ws.Range('I201').Value = 'D'
$ws.Range('J201').Value = 3.75
$ws.Range('K201').Value = 3.4
$ws.Range('L201').Value = 1.8
$ws.Range('M201').Value = 2.8
$ws.Range('N201').Value = 3.1
$ws.Range('O201').Value = 2.25
$ws.Range('P201').Value = 0.25
$ws.Range('Q201').Value = 1.8
$ws.Range('R201').Value = 2
$ws.Range('S201').Value = 2.25
$ws.Range('T201').Value = 1.775
$ws.Range('U201').Value = 2.025
$ws.Range('V201').Value = -1
$ws.Range('W201').Value = 2.1
$ws.Range('X201').Value = -1
$ws.Range('Y201').Value = 0.4
$ws.Range('Z201').Value = -0.5
$ws.Range('AA201').Value = -1
$ws.Range('AB201').Value = 1.025

# Row 268 <= (old) Row 271 data  [id 8162892]
$ws.Range('B268').Value = 8162892
$ws.Range('E268').Value = 'Alajuelense'
$ws.Range('F268').Value = 'AD Guanacasteca'
$ws.Range('G268').Value = 5
$ws.Range('H268').Value = 0
$ws.Range('I268').Value = 'H'
$ws.Range('J268').Value = 1.25
$ws.Range('K268').Value = 5
$ws.Range('L268').Value = 10
$ws.Range('M268').Value = 1.3
$ws.Range('N268').Value = 4.75
$ws.Range('O268').Value = 8
$ws.Range('P268').Value = -1.5
$ws.Range('Q268').Value = 1.9
$ws.Range('R268').Value = 1.9
$ws.Range('S268').Value = 3
$ws.Range('T268').Value = 1.9
$ws.Range('U268').Value = 1.9
$ws.Range('V268').Value = 0.3
$ws.Range('W268').Value = -1
$ws.Range('X268').Value = -1
$ws.Range('Y268').Value = 0.8999999999999999
$ws.Range('Z268').Value = -1
$ws.Range('AA268').Value = 0.8999999999999999
$ws.Range('AB268').Value = -1

# Row 269 <= (old) Row 270 data  [id 8203655]
$ws.Range('B269').Value = 8203655
$ws.Range('E269').Value = 'Municipal Perez Zeledon'
$ws.Range('F269').Value = 'Municipal Liberia'
$ws.Range('G269').Value = 0
$ws.Range('H269').Value = 3
$ws.Range('I269').Value = 'A'
$ws.Range('J269').Value = 3.3
$ws.Range('K269').Value = 3.5
$ws.Range('L269').Value = 2
$ws.Range('M269').Value = 2.9
$ws.Range('N269').Value = 3.4
$ws.Range('O269').Value = 2.2
$ws.Range('P269').Value = 0.25
$ws.Range('Q269').Value = 1.825
$ws.Range('R269').Value = 1.975
$ws.Range('S269').Value = 2.75
$ws.Range('T269').Value = 2
$ws.Range('U269').Value = 1.8
$ws.Range('V269').Value = -1
$ws.Range('W269').Value = -1
$ws.Range('X269').Value = 1.2
$ws.Range('Y269').Value = -1
$ws.Range('Z269').Value = 0.9750000000000001
$ws.Range('AA269').Value = 0.5
$ws.Range('AB269').Value = -0.5

# Row 270 <= (old) Row 268 data  [id 8162893]
$ws.Range('B270').Value = 8162893
$ws.Range('E270').Value = 'AD Grecia'
$ws.Range('F270').Value = 'AD San Carlos'
$ws.Range('G270').Value = 2
$ws.Range('H270').Value = 2
$ws.Range('I270').Value = 'D'
$ws.Range('J270').Value = 5
$ws.Range('K270').Value = 4
$ws.Range('L270').Value = 1.533
$ws.Range('M270').Value = 4.2
$ws.Range('N270').Value = 4.2
$ws.Range('O270').Value = 1.6
$ws.Range('P270').Value = 1
$ws.Range('Q270').Value = 1.775
$ws.Range('R270').Value = 2.025
$ws.Range('S270').Value = 3
$ws.Range('T270').Value = 1.925
$ws.Range('U270').Value = 1.875
$ws.Range('V270').Value = -1
$ws.Range('W270').Value = 3.2
$ws.Range('X270').Value = -1
$ws.Range('Y270').Value = 0.7749999999999999
$ws.Range('Z270').Value = -1
$ws.Range('AA270').Value = 0.925
$ws.Range('AB270').Value = -1

# Row 271 <= (old) Row 269 data  [id 8162895]
$ws.Range('B271').Value = 8162895
$ws.Range('E271').Value = 'Sporting San Jose'
$ws.Range('F271').Value = 'Herediano'
$ws.Range('G271').Value = 1
$ws.Range('H271').Value = 1
$ws.Range('I271').Value = 'D'
$ws.Range('J271').Value = 3.6
$ws.Range('K271').Value = 3.5
$ws.Range('L271').Value = 1.833
$ws.Range('M271').Value = 4.5
$ws.Range('N271').Value = 3.8
$ws.Range('O271').Value = 1.571
$ws.Range('P271').Value = 0.75
$ws.Range('Q271').Value = 2.025
$ws.Range('R271').Value = 1.775
$ws.Range('S271').Value = 2.75
$ws.Range('T271').Value = 1.975
$ws.Range('U271').Value = 1.825
$ws.Range('V271').Value = -1
$ws.Range('W271').Value = 2.8
$ws.Range('X271').Value = -1
$ws.Range('Y271').Value = 1.025
$ws.Range('Z271').Value = -1
$ws.Range('AA271').Value = -1
$ws.Range('AB271').Value = 0.825
